# Updated cryptos list on Tue May  7 02:56:07 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row. Column D values are assigned with a leading apostrophe and then
# restyled to "Normal" so Excel keeps them as text (matching the original
# inlineStr cells) instead of auto-converting numeric-looking strings like
# "592.03" into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.659.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = "'3.084.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = "'592.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").Value = "'156.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.61%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = "'0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").Value = "'3.084.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.85%  '
$ws.Range("E10").Value = '  -3.56%  '
$ws.Range("D11").Value = "'5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("D12").Value = "'0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = "'37.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("E14").Value = '  -3.24%  '
$ws.Range("D15").Value = "'3.599.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("E16").Value = '  -1.80%  '
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").Value = "'63.657.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = "'3.086.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.18%  '
$ws.Range("D20").Value = "'478.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.16%  '
$ws.Range("D21").Value = "'14.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").Value = "'7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  +2.32%  '
$ws.Range("D25").Value = "'81.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").Value = "'12.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("E27").Value = '  +4.70%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = "'7.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.34%  '
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").Value = "'0.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").Value = "'27.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("D35").Value = "'0.0₃0853"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.70%  '
$ws.Range("D36").Value = "'3.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.00%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("D41").Value = "'50.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("D42").Value = "'446.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.11%  '
$ws.Range("D43").Value = "'41.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.53%  '
$ws.Range("E44").Value = '  -1.89%  '
$ws.Range("E45").Value = '  -3.13%  '
$ws.Range("E46").Value = '  +3.82%  '
$ws.Range("D47").Value = "'2.823.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.81%  '
$ws.Range("D48").Value = "'131.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").Value = "'25.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E51").Value = '  +1.17%  '
